# Daily attendance processing - 2025-11-02 01:24:39
# Normalize the "Recorded By" (column G) entries so that "System" is
# listed first among the recorder names/emails for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "backup@backdoor.com, System") {
        $cell.Value = "System, backup@backdoor.com"
    }
    elseif ($val -eq "System, system, backup@backdoor.com") {
        $cell.Value = "System, backup@backdoor.com, system"
    }
}
